$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "90.647.38"
$ws.Range("E2").Value = "  -0.46%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.107.44"
$ws.Range("E3").Value = "  -1.91%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "236.22"
$ws.Range("E5").Value = "  +8.68%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "629.52"
$ws.Range("E6").Value = "  +0.47%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.04"
$ws.Range("E7").Value = "  +1.34%  "
$ws.Range("E8").Value = "  -3.54%  "
$ws.Range("E9").Value = "  +0.11%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "3.107.93"
$ws.Range("E10").Value = "  -1.60%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.713"
$ws.Range("E11").Value = "  -4.35%  "
$ws.Range("E12").Value = "  -1.53%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "36.43"
$ws.Range("E13").Value = "  +4.16%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000244"
$ws.Range("E14").Value = "  -2.56%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.47"
$ws.Range("E15").Value = "  -1.22%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "90.587.07"
$ws.Range("E16").Value = "  -0.36%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.689.23"
$ws.Range("E17").Value = "  -0.72%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.120.03"
$ws.Range("E18").Value = "  -1.30%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.73"
$ws.Range("E19").Value = "  -0.68%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.08"
$ws.Range("E20").Value = "  -1.55%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0000207"
$ws.Range("E21").Value = "  -4.72%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "440.66"
$ws.Range("E22").Value = "  -1.58%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.53"
$ws.Range("E23").Value = "  +5.81%  "
$ws.Range("E24").Value = "  -0.38%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.88"
$ws.Range("E25").Value = "  -4.30%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "87.55"
$ws.Range("E26").Value = "  -1.49%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.34"
$ws.Range("E27").Value = "  -0.18%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.291.70"
$ws.Range("E28").Value = "  -0.26%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  +0.11%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.46"
$ws.Range("E30").Value = "  +3.27%  "
$ws.Range("E31").Value = "  -3.19%  "
$ws.Range("B32").Value = "EthereumClassic"
$ws.Range("C32").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "26.22"
$ws.Range("E32").Value = "  +1.59%  "
$ws.Range("B33").Value = "Stellar"
$ws.Range("C33").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.193"
$ws.Range("E33").Value = "  +24.76%  "
$ws.Range("E34").Value = "  -1.83%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.77"
$ws.Range("E35").Value = "  +0.95%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "507.28"
$ws.Range("E36").Value = "  -4.08%  "
$ws.Range("E37").Value = "  +3.33%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "7.04"
$ws.Range("E38").Value = "  -0.08%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.91"
$ws.Range("E39").Value = "  +0.55%  "
$ws.Range("E40").Value = "  -2.54%  "
$ws.Range("B41").Value = "WhiteBITCoin"
$ws.Range("C41").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "22.17"
$ws.Range("E41").Value = "  -0.40%  "
$ws.Range("B42").Value = "PolygonEcosystemToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.409"
$ws.Range("E42").Value = "  +0.27%  "
$ws.Range("E43").Value = "  -0.04%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0843"
$ws.Range("E44").Value = "  +2.73%  "
$ws.Range("E45").Value = "  +47.20%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.89"
$ws.Range("E46").Value = "  -2.60%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "151.24"
$ws.Range("E47").Value = "  +1.62%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.683"
$ws.Range("E48").Value = "  +5.17%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "45.01"
$ws.Range("E49").Value = "  +1.80%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.33"
$ws.Range("E50").Value = "  -0.57%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.46"
$ws.Range("E51").Value = "  +0.65%  "
